# Update recomputed TPM-derived values for the Ghrl-Gpr39 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.980692333333333
$ws.Range("H2").Value = 5.942077
$ws.Range("I2").Value = 0.05936223023346642
$ws.Range("J2").Value = 0.06064826789587062
$ws.Range("Q2").Value = 0.3696216179387778
$ws.Range("R2").Value = 3.326594561449
$ws.Range("S2").Value = 0.05936223023346642
$ws.Range("T2").Value = 0.06064826789587062

# Row 3
$ws.Range("G3").Value = 2.641089666666666
$ws.Range("H3").Value = 7.923268999999999
$ws.Range("I3").Value = 0.07915463205537174
$ws.Range("J3").Value = 0.0808694570809242
$ws.Range("Q3").Value = 0.4928599052392222
$ws.Range("R3").Value = 4.435739147153
$ws.Range("S3").Value = 0.07915463205537174
$ws.Range("T3").Value = 0.0808694570809242

# Row 4
$ws.Range("G4").Value = 14.555427
$ws.Range("H4").Value = 43.666281
$ws.Range("I4").Value = 0.4362326213815876
$ws.Range("J4").Value = 0.445683270025677
$ws.Range("Q4").Value = 2.716222195133
$ws.Range("R4").Value = 24.445999756197
$ws.Range("S4").Value = 0.4362326213815876
$ws.Range("T4").Value = 0.445683270025677

# Row 5
$ws.Range("G5").Value = 2.1225765
$ws.Range("H5").Value = 4.245153
$ws.Range("I5").Value = 0.06361456181793605
$ws.Range("J5").Value = 0.04332848201108112
$ws.Range("Q5").Value = 0.3960989533435
$ws.Range("R5").Value = 2.376593720061
$ws.Range("S5").Value = 0.06361456181793605
$ws.Range("T5").Value = 0.04332848201108112

# Row 6
$ws.Range("G6").Value = 12.06641933333333
$ws.Range("H6").Value = 36.199258
$ws.Range("I6").Value = 0.3616359545116381
$ws.Range("J6").Value = 0.3694705229864469
$ws.Range("Q6").Value = 2.251742666771778
$ws.Range("R6").Value = 20.265684000946
$ws.Range("S6").Value = 0.3616359545116381
$ws.Range("T6").Value = 0.3694705229864469
